# Update the dSF column (F) values to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -8
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = -7
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = 1
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = 3
